# Actualización automática 2025-07-22 10:10:08
#
# A new client "MILROMER SA" (under advisor GUERRERO FAREZ FABIAN MAURICIO)
# needs to be inserted, alphabetically, right before "MONTESDEOCA ROBLES
# MARIA HILDA" on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets.
# Inserting the row pushes every following row down by one and the running
# totals / "x de N" summary rows move down with it. Because the new client
# has no sales yet, every numeric cell on its row is 0, so the totals in
# the summary rows are unaffected (only the "de 53" -> "de 54" denominator
# text changes, since the roster grew by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (columns A:R, data rows 2-54 -> 2-55,
# summary row 55 -> 56)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(33).Insert()
$ws1.Cells.Item(33, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws1.Cells.Item(33, 2).Value = "MILROMER SA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(33, $c).Value = 0
}

$summary1 = @("3 de 54","0 de 54","4 de 54","1 de 54","0 de 54","2 de 54","3 de 54","0 de 54","0 de 54","2 de 54","14 de 54","0 de 54","0 de 54","1 de 54","0 de 54","0 de 54")
for ($i = 0; $i -lt $summary1.Length; $i++) {
    $ws1.Cells.Item(56, 3 + $i).Value = $summary1[$i]
}

# ---------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (columns A:G, data rows 2-54 -> 2-55,
# totals row 55 -> 56)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(33).Insert()
$ws2.Cells.Item(33, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws2.Cells.Item(33, 2).Value = "MILROMER SA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(33, $c).Value = 0
}

# Sheet 3 ("CUMPLIMIENTO MENSUAL") is unaffected by this update.
